$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# This change re-generates the report: the row that tracked
# f7324243-a653-4e8c-8e9e-f686619523f6 moved up (it is now "In
# Translation") above the rows for 091ce357-... and 1ee9d4d8-...,
# which both remain "Ready for handoff" and shift down by one row.
# This is applied on all three sheets: Overview, zh-cn, de-de.
# -----------------------------------------------------------------

# ---------- Overview sheet ----------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A7").Value2 = "f7324243-a653-4e8c-8e9e-f686619523f6.md"
$ws.Range("B7").Value2 = "In Translation"
$ws.Range("C7").Value2 = "In Translation"

$ws.Range("A8").Value2 = "091ce357-e598-4b5d-aa91-493e68ec530f.md"
$ws.Range("B8").Value2 = "Ready for handoff"
$ws.Range("C8").Value2 = "Ready for handoff"

$ws.Range("A9").Value2 = "1ee9d4d8-34bf-4428-bb8f-8c0209ca05b7.md"
$ws.Range("B9").Value2 = "Ready for handoff"
$ws.Range("C9").Value2 = "Ready for handoff"

# Rebuild the hyperlinks so the displayed text follows the new row
# order while keeping the very same link targets (and therefore the
# same relationship ids) as before, in the same left-to-right,
# top-to-bottom order they originally appeared in.
$targets = @(
  @{Cell="A2"; Addr="https://github.com/OpenLocalizationTest/oltest/blob/0fe490ef8621f44df427f644be1d43fcf77e8f0f/e2e/4938b188-5036-4e4d-9668-f47d24862ac6.md"; Disp="4938b188-5036-4e4d-9668-f47d24862ac6.md"},
  @{Cell="A3"; Addr="https://github.com/OpenLocalizationTest/oltest/blob/b93d6c664216942619d873bff1d775075f5579ae/e2e/373bab7a-ff9b-41bc-86fd-f3488ef6e09f.md"; Disp="373bab7a-ff9b-41bc-86fd-f3488ef6e09f.md"},
  @{Cell="A4"; Addr="https://github.com/OpenLocalizationTest/oltest/blob/4f9590a38b92378f9309447449962cb4753af0da/e2e/51ff0449-1c10-4886-a2bb-0879b1ac49f9.md"; Disp="51ff0449-1c10-4886-a2bb-0879b1ac49f9.md"},
  @{Cell="A5"; Addr="https://github.com/OpenLocalizationTest/oltest/blob/4f9590a38b92378f9309447449962cb4753af0da/e2e/7a5df504-27d3-4f46-8dcd-78ea4380eaa4.md"; Disp="7a5df504-27d3-4f46-8dcd-78ea4380eaa4.md"},
  @{Cell="A6"; Addr="https://github.com/OpenLocalizationTest/oltest/blob/c8bb79716f2d2bfe1582933bbb8f45cc04cf230c/e2e/bf9157c1-8f39-4b76-b675-f98ab404ad41.md"; Disp="bf9157c1-8f39-4b76-b675-f98ab404ad41.md"},
  @{Cell="A7"; Addr="https://github.com/OpenLocalizationTest/oltest/blob/a6d9b5d3ee2d9f7126f53d39dc5eaeaa0f376704/e2e/091ce357-e598-4b5d-aa91-493e68ec530f.md"; Disp="f7324243-a653-4e8c-8e9e-f686619523f6.md"},
  @{Cell="A8"; Addr="https://github.com/OpenLocalizationTest/oltest/blob/d36a81d3c9e28fb4c9170102fa2e007221adf262/e2e/1ee9d4d8-34bf-4428-bb8f-8c0209ca05b7.md"; Disp="091ce357-e598-4b5d-aa91-493e68ec530f.md"},
  @{Cell="A9"; Addr="https://github.com/OpenLocalizationTest/oltest/blob/7e54b9d35931af04890daef61fd43f960c20320f/e2e/f7324243-a653-4e8c-8e9e-f686619523f6.md"; Disp="1ee9d4d8-34bf-4428-bb8f-8c0209ca05b7.md"},
  @{Cell="A10"; Addr="https://github.com/OpenLocalizationTest/oltest/blob/7e54b9d35931af04890daef61fd43f960c20320f/.localization-config"; Disp=".localization-config"}
)

$ws.Hyperlinks.Delete()
foreach ($t in $targets) {
  $ws.Hyperlinks.Add($ws.Range($t.Cell), $t.Addr, "", "", $t.Disp) | Out-Null
}

# ---------- zh-cn sheet ----------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A7").Value2 = "f7324243-a653-4e8c-8e9e-f686619523f6.md"
$ws.Range("B7").Value2 = "In Translation"
$ws.Range("C7").Value2 = "f7324243-a653-4e8c-8e9e-f686619523f6.b0dced80b65e25416406e3268a98ae1b620e7aa4.zh-cn.xlf"
$ws.Range("D7").Value2 = "2016-03-10 04:47:58"

$ws.Range("A8").Value2 = "091ce357-e598-4b5d-aa91-493e68ec530f.md"
$ws.Range("B8").Value2 = "Ready for handoff"
$ws.Range("C8").Value2 = "091ce357-e598-4b5d-aa91-493e68ec530f.af0c841c5b004627d0df61672336ce11aac8f09c.zh-cn.xlf"
$ws.Range("D8").Value2 = "2016-03-10 04:39:23"

$ws.Range("A9").Value2 = "1ee9d4d8-34bf-4428-bb8f-8c0209ca05b7.md"
$ws.Range("B9").Value2 = "Ready for handoff"
$ws.Range("C9").Value2 = "1ee9d4d8-34bf-4428-bb8f-8c0209ca05b7.5e33e09ab16582ff3402d00cc0963bb4d21a49f7.zh-cn.xlf"
$ws.Range("D9").Value2 = "2016-03-10 04:45:42"

$targets = @(
  @{Cell="A2"; Addr="https://github.com/OpenLocalizationTest/oltest/blob/0fe490ef8621f44df427f644be1d43fcf77e8f0f/e2e/4938b188-5036-4e4d-9668-f47d24862ac6.md"; Disp="4938b188-5036-4e4d-9668-f47d24862ac6.md"},
  @{Cell="C2"; Addr="https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a05a60c471731510f7bc322d6b949f86f00a5d25/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/4938b188-5036-4e4d-9668-f47d24862ac6.db788d3874e9e260f3f5a2569d4ee447ace4fa49.zh-cn.xlf"; Disp="4938b188-5036-4e4d-9668-f47d24862ac6.db788d3874e9e260f3f5a2569d4ee447ace4fa49.zh-cn.xlf"},
  @{Cell="E2"; Addr="https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/b7c8cd2eb25be2359ba3b859056bee46b3e8ffe5/e2e/4938b188-5036-4e4d-9668-f47d24862ac6.md"; Disp="4938b188-5036-4e4d-9668-f47d24862ac6.md"},
  @{Cell="F2"; Addr="https://github.com/OpenLocalizationTestOrg/olhandback/blob/06d214cb0cca6f60d9b9656774fd08ee75936331/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/4938b188-5036-4e4d-9668-f47d24862ac6.db788d3874e9e260f3f5a2569d4ee447ace4fa49.zh-cn.xlf"; Disp="4938b188-5036-4e4d-9668-f47d24862ac6.db788d3874e9e260f3f5a2569d4ee447ace4fa49.zh-cn.xlf"},
  @{Cell="A3"; Addr="https://github.com/OpenLocalizationTest/oltest/blob/b93d6c664216942619d873bff1d775075f5579ae/e2e/373bab7a-ff9b-41bc-86fd-f3488ef6e09f.md"; Disp="373bab7a-ff9b-41bc-86fd-f3488ef6e09f.md"},
  @{Cell="C3"; Addr="https://github.com/OpenLocalizationTestOrg/olhandoff/blob/63bf2f4c944e00b3429a3d34462e27e3a4b7ea2c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/373bab7a-ff9b-41bc-86fd-f3488ef6e09f.edbb6d147cb3a96508cbf8f8b684d65ecbb6ad17.zh-cn.xlf"; Disp="373bab7a-ff9b-41bc-86fd-f3488ef6e09f.edbb6d147cb3a96508cbf8f8b684d65ecbb6ad17.zh-cn.xlf"},
  @{Cell="E3"; Addr="https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/3c1bf6c7bb3e8147ec6cfc73792b95a7a27eb0ba/e2e/373bab7a-ff9b-41bc-86fd-f3488ef6e09f.md"; Disp="373bab7a-ff9b-41bc-86fd-f3488ef6e09f.md"},
  @{Cell="F3"; Addr="https://github.com/OpenLocalizationTestOrg/olhandback/blob/e7302fc0ac8e14e8c7450371f4c886b5a9765f20/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/373bab7a-ff9b-41bc-86fd-f3488ef6e09f.edbb6d147cb3a96508cbf8f8b684d65ecbb6ad17.zh-cn.xlf"; Disp="373bab7a-ff9b-41bc-86fd-f3488ef6e09f.edbb6d147cb3a96508cbf8f8b684d65ecbb6ad17.zh-cn.xlf"},
  @{Cell="A4"; Addr="https://github.com/OpenLocalizationTest/oltest/blob/4f9590a38b92378f9309447449962cb4753af0da/e2e/51ff0449-1c10-4886-a2bb-0879b1ac49f9.md"; Disp="51ff0449-1c10-4886-a2bb-0879b1ac49f9.md"},
  @{Cell="C4"; Addr="https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b4ee2c46659cb97e6a2bbba4386746f76739c9ed/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/51ff0449-1c10-4886-a2bb-0879b1ac49f9.6f198a353d0710891a45068f79fc908cb3a9fdc3.zh-cn.xlf"; Disp="51ff0449-1c10-4886-a2bb-0879b1ac49f9.6f198a353d0710891a45068f79fc908cb3a9fdc3.zh-cn.xlf"},
  @{Cell="A5"; Addr="https://github.com/OpenLocalizationTest/oltest/blob/4f9590a38b92378f9309447449962cb4753af0da/e2e/7a5df504-27d3-4f46-8dcd-78ea4380eaa4.md"; Disp="7a5df504-27d3-4f46-8dcd-78ea4380eaa4.md"},
  @{Cell="C5"; Addr="https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b4ee2c46659cb97e6a2bbba4386746f76739c9ed/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/7a5df504-27d3-4f46-8dcd-78ea4380eaa4.73164cee1a14162e06b7cd5a77b87e3171c65b07.zh-cn.xlf"; Disp="7a5df504-27d3-4f46-8dcd-78ea4380eaa4.73164cee1a14162e06b7cd5a77b87e3171c65b07.zh-cn.xlf"},
  @{Cell="A6"; Addr="https://github.com/OpenLocalizationTest/oltest/blob/c8bb79716f2d2bfe1582933bbb8f45cc04cf230c/e2e/bf9157c1-8f39-4b76-b675-f98ab404ad41.md"; Disp="bf9157c1-8f39-4b76-b675-f98ab404ad41.md"},
  @{Cell="C6"; Addr="https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b99bd9fb6603c4029c9a5a13f2f55977db1dc360/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/bf9157c1-8f39-4b76-b675-f98ab404ad41.e578251b7c2a4337bf2f03ef35b3f56d40542222.zh-cn.xlf"; Disp="bf9157c1-8f39-4b76-b675-f98ab404ad41.e578251b7c2a4337bf2f03ef35b3f56d40542222.zh-cn.xlf"},
  @{Cell="A7"; Addr="https://github.com/OpenLocalizationTest/oltest/blob/a6d9b5d3ee2d9f7126f53d39dc5eaeaa0f376704/e2e/091ce357-e598-4b5d-aa91-493e68ec530f.md"; Disp="f7324243-a653-4e8c-8e9e-f686619523f6.md"},
  @{Cell="C7"; Addr="https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d3738be9e50e01023a53bc653be760f54bc8a957/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/091ce357-e598-4b5d-aa91-493e68ec530f.af0c841c5b004627d0df61672336ce11aac8f09c.zh-cn.xlf"; Disp="f7324243-a653-4e8c-8e9e-f686619523f6.b0dced80b65e25416406e3268a98ae1b620e7aa4.zh-cn.xlf"},
  @{Cell="A8"; Addr="https://github.com/OpenLocalizationTest/oltest/blob/d36a81d3c9e28fb4c9170102fa2e007221adf262/e2e/1ee9d4d8-34bf-4428-bb8f-8c0209ca05b7.md"; Disp="091ce357-e598-4b5d-aa91-493e68ec530f.md"},
  @{Cell="C8"; Addr="https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ec4af144d428cd3eda4d67e53f6e4c7dabd6dc6e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/1ee9d4d8-34bf-4428-bb8f-8c0209ca05b7.5e33e09ab16582ff3402d00cc0963bb4d21a49f7.zh-cn.xlf"; Disp="091ce357-e598-4b5d-aa91-493e68ec530f.af0c841c5b004627d0df61672336ce11aac8f09c.zh-cn.xlf"},
  @{Cell="A9"; Addr="https://github.com/OpenLocalizationTest/oltest/blob/7e54b9d35931af04890daef61fd43f960c20320f/e2e/f7324243-a653-4e8c-8e9e-f686619523f6.md"; Disp="1ee9d4d8-34bf-4428-bb8f-8c0209ca05b7.md"},
  @{Cell="C9"; Addr="https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7219fca34662ffce94f0b548741d75bf18400b8a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/f7324243-a653-4e8c-8e9e-f686619523f6.b0dced80b65e25416406e3268a98ae1b620e7aa4.zh-cn.xlf"; Disp="1ee9d4d8-34bf-4428-bb8f-8c0209ca05b7.5e33e09ab16582ff3402d00cc0963bb4d21a49f7.zh-cn.xlf"},
  @{Cell="A10"; Addr="https://github.com/OpenLocalizationTest/oltest/blob/7e54b9d35931af04890daef61fd43f960c20320f/.localization-config"; Disp=".localization-config"}
)

$ws.Hyperlinks.Delete()
foreach ($t in $targets) {
  $ws.Hyperlinks.Add($ws.Range($t.Cell), $t.Addr, "", "", $t.Disp) | Out-Null
}

# ---------- de-de sheet ----------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A7").Value2 = "f7324243-a653-4e8c-8e9e-f686619523f6.md"
$ws.Range("B7").Value2 = "In Translation"
$ws.Range("C7").Value2 = "f7324243-a653-4e8c-8e9e-f686619523f6.b0dced80b65e25416406e3268a98ae1b620e7aa4.de-de.xlf"
$ws.Range("D7").Value2 = "2016-03-10 04:48:07"

$ws.Range("A8").Value2 = "091ce357-e598-4b5d-aa91-493e68ec530f.md"
$ws.Range("B8").Value2 = "Ready for handoff"
$ws.Range("C8").Value2 = "091ce357-e598-4b5d-aa91-493e68ec530f.af0c841c5b004627d0df61672336ce11aac8f09c.de-de.xlf"
$ws.Range("D8").Value2 = "2016-03-10 04:39:31"

$ws.Range("A9").Value2 = "1ee9d4d8-34bf-4428-bb8f-8c0209ca05b7.md"
$ws.Range("B9").Value2 = "Ready for handoff"
$ws.Range("C9").Value2 = "1ee9d4d8-34bf-4428-bb8f-8c0209ca05b7.5e33e09ab16582ff3402d00cc0963bb4d21a49f7.de-de.xlf"
$ws.Range("D9").Value2 = "2016-03-10 04:45:50"

$targets = @(
  @{Cell="A2"; Addr="https://github.com/OpenLocalizationTest/oltest/blob/0fe490ef8621f44df427f644be1d43fcf77e8f0f/e2e/4938b188-5036-4e4d-9668-f47d24862ac6.md"; Disp="4938b188-5036-4e4d-9668-f47d24862ac6.md"},
  @{Cell="C2"; Addr="https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6c1a0d2e85395349d71958495bd5858be675e095/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/4938b188-5036-4e4d-9668-f47d24862ac6.db788d3874e9e260f3f5a2569d4ee447ace4fa49.de-de.xlf"; Disp="4938b188-5036-4e4d-9668-f47d24862ac6.db788d3874e9e260f3f5a2569d4ee447ace4fa49.de-de.xlf"},
  @{Cell="E2"; Addr="https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/1e573a59995dbffc0f34db50b030225da385e6d4/e2e/4938b188-5036-4e4d-9668-f47d24862ac6.md"; Disp="4938b188-5036-4e4d-9668-f47d24862ac6.md"},
  @{Cell="F2"; Addr="https://github.com/OpenLocalizationTestOrg/olhandback/blob/a55205657936bf06827d5161ae943eb28a3da1af/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/4938b188-5036-4e4d-9668-f47d24862ac6.db788d3874e9e260f3f5a2569d4ee447ace4fa49.de-de.xlf"; Disp="4938b188-5036-4e4d-9668-f47d24862ac6.db788d3874e9e260f3f5a2569d4ee447ace4fa49.de-de.xlf"},
  @{Cell="A3"; Addr="https://github.com/OpenLocalizationTest/oltest/blob/b93d6c664216942619d873bff1d775075f5579ae/e2e/373bab7a-ff9b-41bc-86fd-f3488ef6e09f.md"; Disp="373bab7a-ff9b-41bc-86fd-f3488ef6e09f.md"},
  @{Cell="C3"; Addr="https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c42d6853839fd15a9eb59ec28a3a7b8196ae468c/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/373bab7a-ff9b-41bc-86fd-f3488ef6e09f.edbb6d147cb3a96508cbf8f8b684d65ecbb6ad17.de-de.xlf"; Disp="373bab7a-ff9b-41bc-86fd-f3488ef6e09f.edbb6d147cb3a96508cbf8f8b684d65ecbb6ad17.de-de.xlf"},
  @{Cell="E3"; Addr="https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/d4a7ef6fca95bc477edecd1e6207541212d128b5/e2e/373bab7a-ff9b-41bc-86fd-f3488ef6e09f.md"; Disp="373bab7a-ff9b-41bc-86fd-f3488ef6e09f.md"},
  @{Cell="F3"; Addr="https://github.com/OpenLocalizationTestOrg/olhandback/blob/20c2391ce86313c2a460ac890ab87c1578235033/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/373bab7a-ff9b-41bc-86fd-f3488ef6e09f.edbb6d147cb3a96508cbf8f8b684d65ecbb6ad17.de-de.xlf"; Disp="373bab7a-ff9b-41bc-86fd-f3488ef6e09f.edbb6d147cb3a96508cbf8f8b684d65ecbb6ad17.de-de.xlf"},
  @{Cell="A4"; Addr="https://github.com/OpenLocalizationTest/oltest/blob/4f9590a38b92378f9309447449962cb4753af0da/e2e/51ff0449-1c10-4886-a2bb-0879b1ac49f9.md"; Disp="51ff0449-1c10-4886-a2bb-0879b1ac49f9.md"},
  @{Cell="C4"; Addr="https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8e6ea7fc2ef40a502d09c27fd3fa1a628d50c6eb/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/51ff0449-1c10-4886-a2bb-0879b1ac49f9.6f198a353d0710891a45068f79fc908cb3a9fdc3.de-de.xlf"; Disp="51ff0449-1c10-4886-a2bb-0879b1ac49f9.6f198a353d0710891a45068f79fc908cb3a9fdc3.de-de.xlf"},
  @{Cell="A5"; Addr="https://github.com/OpenLocalizationTest/oltest/blob/4f9590a38b92378f9309447449962cb4753af0da/e2e/7a5df504-27d3-4f46-8dcd-78ea4380eaa4.md"; Disp="7a5df504-27d3-4f46-8dcd-78ea4380eaa4.md"},
  @{Cell="C5"; Addr="https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8e6ea7fc2ef40a502d09c27fd3fa1a628d50c6eb/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/7a5df504-27d3-4f46-8dcd-78ea4380eaa4.73164cee1a14162e06b7cd5a77b87e3171c65b07.de-de.xlf"; Disp="7a5df504-27d3-4f46-8dcd-78ea4380eaa4.73164cee1a14162e06b7cd5a77b87e3171c65b07.de-de.xlf"},
  @{Cell="A6"; Addr="https://github.com/OpenLocalizationTest/oltest/blob/c8bb79716f2d2bfe1582933bbb8f45cc04cf230c/e2e/bf9157c1-8f39-4b76-b675-f98ab404ad41.md"; Disp="bf9157c1-8f39-4b76-b675-f98ab404ad41.md"},
  @{Cell="C6"; Addr="https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fcfe79993655e1f007d76fc4d7783775e58464a3/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/bf9157c1-8f39-4b76-b675-f98ab404ad41.e578251b7c2a4337bf2f03ef35b3f56d40542222.de-de.xlf"; Disp="bf9157c1-8f39-4b76-b675-f98ab404ad41.e578251b7c2a4337bf2f03ef35b3f56d40542222.de-de.xlf"},
  @{Cell="A7"; Addr="https://github.com/OpenLocalizationTest/oltest/blob/a6d9b5d3ee2d9f7126f53d39dc5eaeaa0f376704/e2e/091ce357-e598-4b5d-aa91-493e68ec530f.md"; Disp="f7324243-a653-4e8c-8e9e-f686619523f6.md"},
  @{Cell="C7"; Addr="https://github.com/OpenLocalizationTestOrg/olhandoff/blob/691e331c0beeeb65bf4d8495475b42a3279530db/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/091ce357-e598-4b5d-aa91-493e68ec530f.af0c841c5b004627d0df61672336ce11aac8f09c.de-de.xlf"; Disp="f7324243-a653-4e8c-8e9e-f686619523f6.b0dced80b65e25416406e3268a98ae1b620e7aa4.de-de.xlf"},
  @{Cell="A8"; Addr="https://github.com/OpenLocalizationTest/oltest/blob/d36a81d3c9e28fb4c9170102fa2e007221adf262/e2e/1ee9d4d8-34bf-4428-bb8f-8c0209ca05b7.md"; Disp="091ce357-e598-4b5d-aa91-493e68ec530f.md"},
  @{Cell="C8"; Addr="https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1433806b258c4f0bd8bf365b81e9d9286d159a58/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/1ee9d4d8-34bf-4428-bb8f-8c0209ca05b7.5e33e09ab16582ff3402d00cc0963bb4d21a49f7.de-de.xlf"; Disp="091ce357-e598-4b5d-aa91-493e68ec530f.af0c841c5b004627d0df61672336ce11aac8f09c.de-de.xlf"},
  @{Cell="A9"; Addr="https://github.com/OpenLocalizationTest/oltest/blob/7e54b9d35931af04890daef61fd43f960c20320f/e2e/f7324243-a653-4e8c-8e9e-f686619523f6.md"; Disp="1ee9d4d8-34bf-4428-bb8f-8c0209ca05b7.md"},
  @{Cell="C9"; Addr="https://github.com/OpenLocalizationTestOrg/olhandoff/blob/92c98aec1bb18ea4965cfe3b170816d20487acb7/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/f7324243-a653-4e8c-8e9e-f686619523f6.b0dced80b65e25416406e3268a98ae1b620e7aa4.de-de.xlf"; Disp="1ee9d4d8-34bf-4428-bb8f-8c0209ca05b7.5e33e09ab16582ff3402d00cc0963bb4d21a49f7.de-de.xlf"},
  @{Cell="A10"; Addr="https://github.com/OpenLocalizationTest/oltest/blob/7e54b9d35931af04890daef61fd43f960c20320f/.localization-config"; Disp=".localization-config"}
)

$ws.Hyperlinks.Delete()
foreach ($t in $targets) {
  $ws.Hyperlinks.Add($ws.Range($t.Cell), $t.Addr, "", "", $t.Disp) | Out-Null
}

$wb.Save()
